$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 498.74545
$ws.Cells.Item(17, 10).Value = 498.74545
$ws.Cells.Item(17, 12).Value = 1496.23635
$ws.Cells.Item(17, 14).Value = -1832.23635
$ws.Cells.Item(116, 8).Value = 2063.182
$ws.Cells.Item(116, 9).Value = 1764.6666
$ws.Cells.Item(116, 11).Value = 1764.6666
$ws.Cells.Item(116, 13).Value = 1677.3334
$ws.Cells.Item(132, 8).Value = 876774.75
$ws.Cells.Item(132, 9).Value = 1626.5883
$ws.Cells.Item(132, 10).Value = 9803286
$ws.Cells.Item(132, 11).Value = 4879.7649
$ws.Cells.Item(132, 12).Value = 29409858
$ws.Cells.Item(132, 13).Value = -2349.7649
$ws.Cells.Item(132, 14).Value = -29414918
$ws.Cells.Item(137, 8).Value = 2501595.8
$ws.Cells.Item(137, 9).Value = 4349092
$ws.Cells.Item(137, 10).Value = 2041.9412
$ws.Cells.Item(137, 11).Value = 13047276
$ws.Cells.Item(137, 12).Value = 6125.8236
$ws.Cells.Item(137, 13).Value = -13044726
$ws.Cells.Item(137, 14).Value = -11225.8236
$ws.Cells.Item(138, 8).Value = 3088249.5
$ws.Cells.Item(138, 9).Value = 1371.5428
$ws.Cells.Item(138, 10).Value = 8774604
$ws.Cells.Item(138, 11).Value = 4114.6284
$ws.Cells.Item(138, 12).Value = 26323812
$ws.Cells.Item(138, 13).Value = 1025.3716
$ws.Cells.Item(138, 14).Value = -26334092
$ws.Cells.Item(141, 8).Value = 455.13727
$ws.Cells.Item(141, 9).Value = 421.34784
$ws.Cells.Item(141, 10).Value = 766
$ws.Cells.Item(141, 11).Value = 1264.04352
$ws.Cells.Item(141, 12).Value = 2298
$ws.Cells.Item(141, 13).Value = 3915.95648
$ws.Cells.Item(141, 14).Value = -12658

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 911.27
$ws.Cells.Item(32, 9).Value = 796.3253
$ws.Cells.Item(32, 10).Value = 1472.4706
$ws.Cells.Item(32, 11).Value = 796.3253
$ws.Cells.Item(32, 12).Value = 1472.4706
$ws.Cells.Item(32, 13).Value = -509.3253
$ws.Cells.Item(32, 14).Value = -2046.4706
$ws.Cells.Item(45, 8).Value = 1398.5294
$ws.Cells.Item(45, 9).Value = 931.25
$ws.Cells.Item(45, 10).Value = 2520
$ws.Cells.Item(45, 11).Value = 931.25
$ws.Cells.Item(45, 12).Value = 2520
$ws.Cells.Item(45, 13).Value = -554.25
$ws.Cells.Item(45, 14).Value = -3274
$ws.Cells.Item(74, 8).Value = 7876747.5
$ws.Cells.Item(74, 9).Value = 11410591
$ws.Cells.Item(74, 10).Value = 102291.4
$ws.Cells.Item(74, 11).Value = 11410591
$ws.Cells.Item(74, 12).Value = 102291.4
$ws.Cells.Item(74, 13).Value = -11409717
$ws.Cells.Item(74, 14).Value = -104039.4
$ws.Cells.Item(77, 8).Value = 7876747.5
$ws.Cells.Item(77, 9).Value = 11410591
$ws.Cells.Item(77, 10).Value = 102291.4
$ws.Cells.Item(77, 11).Value = 57052955
$ws.Cells.Item(77, 12).Value = 511457
$ws.Cells.Item(77, 13).Value = -57048587
$ws.Cells.Item(77, 14).Value = -520193
$ws.Cells.Item(122, 8).Value = 3970243.8
$ws.Cells.Item(122, 9).Value = 2039.6086
$ws.Cells.Item(122, 11).Value = 6118.825800000001
$ws.Cells.Item(122, 13).Value = -3668.825800000001
$ws.Cells.Item(132, 8).Value = 78609.19
$ws.Cells.Item(132, 9).Value = 59808.766
$ws.Cells.Item(132, 10).Value = 110569.9
$ws.Cells.Item(132, 11).Value = 179426.298
$ws.Cells.Item(132, 12).Value = 331709.7
$ws.Cells.Item(132, 13).Value = -176896.298
$ws.Cells.Item(132, 14).Value = -336769.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2314.5173
$ws.Cells.Item(134, 9).Value = 1346.5294
$ws.Cells.Item(134, 11).Value = 4039.5882
$ws.Cells.Item(134, 13).Value = -1504.5882

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1647.9014
$ws.Cells.Item(31, 9).Value = 1001.3137
$ws.Cells.Item(31, 10).Value = 3296.7
$ws.Cells.Item(31, 11).Value = 1001.3137
$ws.Cells.Item(31, 12).Value = 3296.7
$ws.Cells.Item(31, 13).Value = -706.3137
$ws.Cells.Item(31, 14).Value = -3886.7
$ws.Cells.Item(34, 8).Value = 1647.9014
$ws.Cells.Item(34, 9).Value = 1001.3137
$ws.Cells.Item(34, 10).Value = 3296.7
$ws.Cells.Item(34, 11).Value = 1001.3137
$ws.Cells.Item(34, 12).Value = 3296.7
$ws.Cells.Item(34, 13).Value = -799.3137
$ws.Cells.Item(34, 14).Value = -3700.7
$ws.Cells.Item(58, 8).Value = 19609520
$ws.Cells.Item(58, 9).Value = 26317366
$ws.Cells.Item(58, 10).Value = 1972.4615
$ws.Cells.Item(58, 11).Value = 26317366
$ws.Cells.Item(58, 12).Value = 1972.4615
$ws.Cells.Item(58, 13).Value = -26317163
$ws.Cells.Item(58, 14).Value = -2378.4615
$ws.Cells.Item(99, 8).Value = 5881.769
$ws.Cells.Item(99, 9).Value = 6411.647
$ws.Cells.Item(99, 10).Value = 4880.8887
$ws.Cells.Item(99, 11).Value = 6411.647
$ws.Cells.Item(99, 12).Value = 4880.8887
$ws.Cells.Item(99, 13).Value = -4913.647
$ws.Cells.Item(99, 14).Value = -7876.8887
$ws.Cells.Item(126, 8).Value = 5881.769
$ws.Cells.Item(126, 9).Value = 6411.647
$ws.Cells.Item(126, 10).Value = 4880.8887
$ws.Cells.Item(126, 11).Value = 19234.941
$ws.Cells.Item(126, 12).Value = 14642.6661
$ws.Cells.Item(126, 13).Value = -16764.941
$ws.Cells.Item(126, 14).Value = -19582.6661
$ws.Cells.Item(132, 8).Value = 43372.457
$ws.Cells.Item(132, 9).Value = 29380.75
$ws.Cells.Item(132, 10).Value = 85347.586
$ws.Cells.Item(132, 11).Value = 88142.25
$ws.Cells.Item(132, 12).Value = 256042.758
$ws.Cells.Item(132, 13).Value = -85612.25
$ws.Cells.Item(132, 14).Value = -261102.758
$ws.Cells.Item(134, 8).Value = 28420.025
$ws.Cells.Item(134, 9).Value = 1409.1212
$ws.Cells.Item(134, 10).Value = 155757.14
$ws.Cells.Item(134, 11).Value = 4227.363600000001
$ws.Cells.Item(134, 12).Value = 467271.42
$ws.Cells.Item(134, 13).Value = -1692.363600000001
$ws.Cells.Item(134, 14).Value = -472341.42
$ws.Cells.Item(136, 8).Value = 19609520
$ws.Cells.Item(136, 9).Value = 26317366
$ws.Cells.Item(136, 10).Value = 1972.4615
$ws.Cells.Item(136, 11).Value = 78952098
$ws.Cells.Item(136, 12).Value = 5917.3845
$ws.Cells.Item(136, 13).Value = -78949548
$ws.Cells.Item(136, 14).Value = -11017.3845

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1145.1111
$ws.Cells.Item(34, 9).Value = 383.33334
$ws.Cells.Item(34, 10).Value = 2668.6667
$ws.Cells.Item(34, 11).Value = 1150.00002
$ws.Cells.Item(34, 12).Value = 8006.000100000001
$ws.Cells.Item(34, 13).Value = -1066.00002
$ws.Cells.Item(34, 14).Value = -8174.000100000001
$ws.Cells.Item(39, 8).Value = 1800
$ws.Cells.Item(39, 10).Value = 2450
$ws.Cells.Item(39, 12).Value = 7350
$ws.Cells.Item(39, 14).Value = -7938
$ws.Cells.Item(55, 8).Value = 1500
$ws.Cells.Item(55, 9).Value = 375
$ws.Cells.Item(55, 10).Value = 2400
$ws.Cells.Item(55, 11).Value = 1125
$ws.Cells.Item(55, 12).Value = 7200
$ws.Cells.Item(55, 13).Value = -948
$ws.Cells.Item(55, 14).Value = -7554
$ws.Cells.Item(131, 8).Value = 1022.5614
$ws.Cells.Item(131, 10).Value = 1218.3256
$ws.Cells.Item(131, 12).Value = 3654.976799999999
$ws.Cells.Item(131, 14).Value = -13734.9768

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 8533.333000000001
$ws.Cells.Item(7, 9).Value = 600
$ws.Cells.Item(7, 11).Value = 600
$ws.Cells.Item(7, 13).Value = -488
$ws.Cells.Item(8, 8).Value = 8533.333000000001
$ws.Cells.Item(8, 9).Value = 600
$ws.Cells.Item(8, 11).Value = 600
$ws.Cells.Item(8, 13).Value = -461
$ws.Cells.Item(13, 8).Value = 444
$ws.Cells.Item(13, 9).Value = 404.44446
$ws.Cells.Item(13, 10).Value = 800
$ws.Cells.Item(13, 11).Value = 404.44446
$ws.Cells.Item(13, 12).Value = 800
$ws.Cells.Item(13, 13).Value = -265.44446
$ws.Cells.Item(13, 14).Value = -1078
$ws.Cells.Item(70, 8).Value = 5455.1377
$ws.Cells.Item(70, 9).Value = 5244.4443
$ws.Cells.Item(70, 10).Value = 5799.909
$ws.Cells.Item(70, 11).Value = 5244.4443
$ws.Cells.Item(70, 12).Value = 5799.909
$ws.Cells.Item(70, 13).Value = -4974.4443
$ws.Cells.Item(70, 14).Value = -6339.909
$ws.Cells.Item(73, 8).Value = 5455.1377
$ws.Cells.Item(73, 9).Value = 5244.4443
$ws.Cells.Item(73, 10).Value = 5799.909
$ws.Cells.Item(73, 11).Value = 5244.4443
$ws.Cells.Item(73, 12).Value = 5799.909
$ws.Cells.Item(73, 13).Value = -4308.4443
$ws.Cells.Item(73, 14).Value = -7671.909
$ws.Cells.Item(122, 8).Value = 3955.3333
$ws.Cells.Item(122, 9).Value = 3509.5293
$ws.Cells.Item(122, 10).Value = 5850
$ws.Cells.Item(122, 11).Value = 10528.5879
$ws.Cells.Item(122, 12).Value = 17550
$ws.Cells.Item(122, 13).Value = -8078.5879
$ws.Cells.Item(122, 14).Value = -22450
$ws.Cells.Item(126, 8).Value = 2182.8
$ws.Cells.Item(126, 9).Value = 1337.3334
$ws.Cells.Item(126, 10).Value = 2746.4443
$ws.Cells.Item(126, 11).Value = 4012.0002
$ws.Cells.Item(126, 12).Value = 8239.332900000001
$ws.Cells.Item(126, 13).Value = -1542.0002
$ws.Cells.Item(126, 14).Value = -13179.3329
$ws.Cells.Item(132, 8).Value = 57270.555
$ws.Cells.Item(132, 9).Value = 43093.418
$ws.Cells.Item(132, 10).Value = 85624.836
$ws.Cells.Item(132, 11).Value = 129280.254
$ws.Cells.Item(132, 12).Value = 256874.508
$ws.Cells.Item(132, 13).Value = -126750.254
$ws.Cells.Item(132, 14).Value = -261934.508

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1788.1072
$ws.Cells.Item(7, 9).Value = 1651.6957
$ws.Cells.Item(7, 10).Value = 2415.6
$ws.Cells.Item(7, 11).Value = 1651.6957
$ws.Cells.Item(7, 12).Value = 2415.6
$ws.Cells.Item(7, 13).Value = -1539.6957
$ws.Cells.Item(7, 14).Value = -2639.6
$ws.Cells.Item(22, 8).Value = 1800.6666
$ws.Cells.Item(22, 9).Value = 900
$ws.Cells.Item(22, 10).Value = 2251
$ws.Cells.Item(22, 11).Value = 900
$ws.Cells.Item(22, 12).Value = 2251
$ws.Cells.Item(22, 13).Value = -605
$ws.Cells.Item(22, 14).Value = -2841
$ws.Cells.Item(27, 8).Value = 1800.6666
$ws.Cells.Item(27, 9).Value = 900
$ws.Cells.Item(27, 10).Value = 2251
$ws.Cells.Item(27, 11).Value = 900
$ws.Cells.Item(27, 12).Value = 2251
$ws.Cells.Item(27, 13).Value = -793
$ws.Cells.Item(27, 14).Value = -2465
$ws.Cells.Item(126, 8).Value = 1788.1072
$ws.Cells.Item(126, 9).Value = 1651.6957
$ws.Cells.Item(126, 10).Value = 2415.6
$ws.Cells.Item(126, 11).Value = 4955.0871
$ws.Cells.Item(126, 12).Value = 7246.799999999999
$ws.Cells.Item(126, 13).Value = -2485.0871
$ws.Cells.Item(126, 14).Value = -12186.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 12356.857
$ws.Cells.Item(14, 10).Value = 9416.5
$ws.Cells.Item(14, 12).Value = 9416.5
$ws.Cells.Item(14, 14).Value = -9752.5
$ws.Cells.Item(132, 8).Value = 53160.543
$ws.Cells.Item(132, 9).Value = 44023.74
$ws.Cells.Item(132, 10).Value = 85490.766
$ws.Cells.Item(132, 11).Value = 132071.22
$ws.Cells.Item(132, 12).Value = 256472.298
$ws.Cells.Item(132, 13).Value = -129541.22
$ws.Cells.Item(132, 14).Value = -261532.298
$ws.Cells.Item(136, 8).Value = 37611.418
$ws.Cells.Item(136, 9).Value = 25737.8
$ws.Cells.Item(136, 10).Value = 69274.39999999999
$ws.Cells.Item(136, 11).Value = 77213.39999999999
$ws.Cells.Item(136, 12).Value = 207823.2
$ws.Cells.Item(136, 13).Value = -74663.39999999999
$ws.Cells.Item(136, 14).Value = -212923.2
$ws.Cells.Item(140, 8).Value = 40238.5
$ws.Cells.Item(140, 10).Value = 40238.5
$ws.Cells.Item(140, 12).Value = 40238.5
$ws.Cells.Item(140, 14).Value = -50598.5
